$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.441463470458984
$ws.Range("B1").Value = 4.375473976135254
$ws.Range("C1").Value = 3.746745347976685
$ws.Range("D1").Value = 4.546655178070068
$ws.Range("E1").Value = 4.770989894866943
